# care_organization_experience.pptx - "final re-run including fig: unmet social support"
#
# All geometry in the source OOXML is expressed in EMU (914400 EMU = 1 inch,
# 12700 EMU = 1 point). The PowerPoint object model stores Left/Top/Width/
# Height in points as a (single-precision) float, and converting back to EMU
# on save truncates rather than rounds - so an exact EMU/12700 division can
# land one EMU short after the round trip. A tiny (sub-EMU) epsilon added to
# every assigned point value keeps the float on the correct side of the
# truncation boundary without being large enough to ever push the result
# into the next EMU.

$EMU    = 12700
$NUDGE  = 0.00002   # << 1/12700 pt; corrects float32 truncation, never overshoots

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(1)
$top = $s.Shapes.Item(1)      # the single top-level group shape on the slide
$grp = $top.GroupItems

function EmuToPt($emu) {
    return ($emu / $EMU) + $NUDGE
}

function Set-TopHeight($name, $top_emu, $cy_emu) {
    $sh = $grp.Item($name)
    if ($top_emu -ne $null) { $sh.Top = EmuToPt $top_emu }
    $sh.Height = EmuToPt $cy_emu
}

function Set-Text($name, $text) {
    $sh = $grp.Item($name)
    $sh.TextFrame.TextRange.Text = $text
}

function Set-Top($name, $top_emu) {
    $sh = $grp.Item($name)
    $sh.Top = EmuToPt $top_emu
}

# ---------------------------------------------------------------------------
# Left-hand column of bars (rc4..rc9) - reposition/resize
# ---------------------------------------------------------------------------
Set-TopHeight "rc4" $null      158410
Set-TopHeight "rc5" 2105187    760370
Set-TopHeight "rc6" 2865558    1679151
Set-TopHeight "rc7" 4544709    1045509
Set-TopHeight "rc8" 5590219    760370
Set-TopHeight "rc9" 6350590    158410

# ---------------------------------------------------------------------------
# Right-hand column of bars (rc10..rc15) - reposition/resize
# ---------------------------------------------------------------------------
Set-TopHeight "rc10" $null     126728
Set-TopHeight "rc11" 2073505   2091019
Set-TopHeight "rc12" 4164524   1805880
Set-TopHeight "rc13" 5970404   221774
Set-TopHeight "rc14" 6192179   190092
Set-TopHeight "rc15" 6382272   126728

# ---------------------------------------------------------------------------
# Left column data-labels (tx16..tx21) - reposition, some text updated
# ---------------------------------------------------------------------------
Set-Top  "tx16" 1942067
# tx16 text unchanged: "3 (2%)"

Set-Top  "tx17" 2401457
Set-Text "tx17" "24 (17%)"

Set-Top  "tx18" 3621218
Set-Text "tx18" "53 (37%)"

Set-Top  "tx19" 4983549
Set-Text "tx19" "33 (23%)"

Set-Top  "tx20" 5886489
Set-Text "tx20" "24 (17%)"

Set-Top  "tx21" 6345880
# tx21 text unchanged: "5 (3%)"

# ---------------------------------------------------------------------------
# Right column data-labels (tx22..tx27) - reposition, some text updated
# ---------------------------------------------------------------------------
Set-Top  "tx22" 1926226
Set-Text "tx22" "4 (3%)"

Set-Top  "tx23" 3035099
Set-Text "tx23" "66 (46%)"

Set-Top  "tx24" 4983549
Set-Text "tx24" "57 (40%)"

Set-Top  "tx25" 5997377
Set-Text "tx25" "7 (5%)"

Set-Top  "tx26" 6203311
# tx26 text unchanged: "6 (4%)"

Set-Top  "tx27" 6361721
# tx27 text unchanged: "4 (3%)"

# ---------------------------------------------------------------------------
# Caption text: sample size changed N=146 -> N=144
# ---------------------------------------------------------------------------
Set-Text "tx51" "healthcare personnel (N=144)"
